# Update the timestamp column (Z) with the new run's timestamps.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamps = @{
    2  = "2025-11-03T00:10:16.536319"
    3  = "2025-11-03T00:10:16.536319"
    4  = "2025-11-03T00:10:16.536319"
    5  = "2025-11-03T00:10:16.536319"
    6  = "2025-11-03T00:10:16.536319"
    7  = "2025-11-03T00:10:16.536319"
    8  = "2025-11-03T00:10:16.536319"
    9  = "2025-11-03T00:10:16.536319"
    10 = "2025-11-03T00:10:16.537319"
    11 = "2025-11-03T00:10:16.537319"
    12 = "2025-11-03T00:10:16.537319"
    13 = "2025-11-03T00:10:16.537319"
    14 = "2025-11-03T00:10:16.537319"
    15 = "2025-11-03T00:10:16.537319"
    16 = "2025-11-03T00:10:16.537319"
    17 = "2025-11-03T00:10:16.537319"
    18 = "2025-11-03T00:10:16.537319"
    19 = "2025-11-03T00:10:16.537319"
    20 = "2025-11-03T00:10:16.538319"
    21 = "2025-11-03T00:10:16.538319"
    22 = "2025-11-03T00:10:16.538319"
    23 = "2025-11-03T00:10:16.538319"
    24 = "2025-11-03T00:10:16.538319"
    25 = "2025-11-03T00:10:16.538319"
    26 = "2025-11-03T00:10:16.538319"
    27 = "2025-11-03T00:10:16.538319"
    28 = "2025-11-03T00:10:16.538319"
    29 = "2025-11-03T00:10:16.538319"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
